$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q4" and push the existing
#    "2021-Q2" row down to row 3 (copy-down preserves the row's styling).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("B2").Copy($summary.Range("B3"))
$summary.Range("C2").Copy($summary.Range("C3"))
$summary.Range("D2").Copy($summary.Range("D3"))

$summary.Range("A3").Value2 = 1
$summary.Range("B2").Value2 = "2022-Q4"
$summary.Range("C2").Value2 = 11
$summary.Range("D2").Value2 = 0.38

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right after "总计" (i.e. before the
#    existing "2021-Q2" sheet, which is the active sheet) and fill it with
#    the fund holdings detail table.
# ---------------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $summary)
$new.Name = "2022-Q4"

# Reuse the "总计" header/index-column formatting (style id 2) for the new
# sheet's header row and index column, same as the rest of the workbook.
$summary.Range("B1").Copy($new.Range("B1:H1"))
$summary.Range("A2").Copy($new.Range("A2:A12"))

# Columns B:G hold text values (fund code, name, size, position %, etc.) -
# force text formatting first so numeric-looking strings (e.g. "009668",
# "0.20") keep their original formatting instead of being parsed as numbers.
$new.Range("B2:G12").NumberFormat = "@"

$new.Range("B1").Value2 = "基金代码"
$new.Range("C1").Value2 = "基金名称"
$new.Range("D1").Value2 = "基金规模"
$new.Range("E1").Value2 = "股票总仓位"
$new.Range("F1").Value2 = "仓位占比"
$new.Range("G1").Value2 = "持有市值(亿元)"
$new.Range("H1").Value2 = "仓位排名"

$new.Range("A2").Value2 = 0
$new.Range("B2").Value2 = "160642"
$new.Range("C2").Value2 = "鹏华增瑞灵活配置混合（LOF）"
$new.Range("D2").Value2 = "2.22"
$new.Range("E2").Value2 = "85.42"
$new.Range("F2").Value2 = "3.24"
$new.Range("G2").Value2 = "0.0719"
$new.Range("H2").Value2 = 8

$new.Range("A3").Value2 = 1
$new.Range("B3").Value2 = "005416"
$new.Range("C3").Value2 = "鹏华尊惠18个月定期开放混合A"
$new.Range("D3").Value2 = "2.83"
$new.Range("E3").Value2 = "39.19"
$new.Range("F3").Value2 = "2.46"
$new.Range("G3").Value2 = "0.0696"
$new.Range("H3").Value2 = 1

$new.Range("A4").Value2 = 2
$new.Range("B4").Value2 = "015026"
$new.Range("C4").Value2 = "鹏华增华混合A"
$new.Range("D4").Value2 = "1.65"
$new.Range("E4").Value2 = "78.57"
$new.Range("F4").Value2 = "3.67"
$new.Range("G4").Value2 = "0.0606"
$new.Range("H4").Value2 = 5

$new.Range("A5").Value2 = 3
$new.Range("B5").Value2 = "009668"
$new.Range("C5").Value2 = "鹏华安庆混合C"
$new.Range("D5").Value2 = "2.79"
$new.Range("E5").Value2 = "39.85"
$new.Range("F5").Value2 = "1.52"
$new.Range("G5").Value2 = "0.0424"
$new.Range("H5").Value2 = 8

$new.Range("A6").Value2 = 4
$new.Range("B6").Value2 = "009667"
$new.Range("C6").Value2 = "鹏华安庆混合A"
$new.Range("D6").Value2 = "2.34"
$new.Range("E6").Value2 = "39.85"
$new.Range("F6").Value2 = "1.52"
$new.Range("G6").Value2 = "0.0356"
$new.Range("H6").Value2 = 8

$new.Range("A7").Value2 = 5
$new.Range("B7").Value2 = "011573"
$new.Range("C7").Value2 = "鹏华安荣混合C"
$new.Range("D7").Value2 = "1.87"
$new.Range("E7").Value2 = "39.61"
$new.Range("F7").Value2 = "1.42"
$new.Range("G7").Value2 = "0.0266"
$new.Range("H7").Value2 = 9

$new.Range("A8").Value2 = 6
$new.Range("B8").Value2 = "011572"
$new.Range("C8").Value2 = "鹏华安荣混合A"
$new.Range("D8").Value2 = "1.54"
$new.Range("E8").Value2 = "39.61"
$new.Range("F8").Value2 = "1.42"
$new.Range("G8").Value2 = "0.0219"
$new.Range("H8").Value2 = 9

$new.Range("A9").Value2 = 7
$new.Range("B9").Value2 = "009231"
$new.Range("C9").Value2 = "鹏华安和混合C"
$new.Range("D9").Value2 = "1.54"
$new.Range("E9").Value2 = "38.20"
$new.Range("F9").Value2 = "1.37"
$new.Range("G9").Value2 = "0.0211"
$new.Range("H9").Value2 = 10

$new.Range("A10").Value2 = 8
$new.Range("B10").Value2 = "009230"
$new.Range("C10").Value2 = "鹏华安和混合A"
$new.Range("D10").Value2 = "1.37"
$new.Range("E10").Value2 = "38.20"
$new.Range("F10").Value2 = "1.37"
$new.Range("G10").Value2 = "0.0188"
$new.Range("H10").Value2 = 10

$new.Range("A11").Value2 = 9
$new.Range("B11").Value2 = "005417"
$new.Range("C11").Value2 = "鹏华尊惠18个月定期开放混合C"
$new.Range("D11").Value2 = "0.36"
$new.Range("E11").Value2 = "39.19"
$new.Range("F11").Value2 = "2.46"
$new.Range("G11").Value2 = "0.0089"
$new.Range("H11").Value2 = 1

$new.Range("A12").Value2 = 10
$new.Range("B12").Value2 = "015027"
$new.Range("C12").Value2 = "鹏华增华混合C"
$new.Range("D12").Value2 = "0.20"
$new.Range("E12").Value2 = "78.57"
$new.Range("F12").Value2 = "3.67"
$new.Range("G12").Value2 = "0.0073"
$new.Range("H12").Value2 = 5

# ---------------------------------------------------------------------------
# 3) "2021-Q2" was the selected tab before the edit - keep it selected
#    (adding the new sheet would otherwise steal the active-sheet focus).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
